# Release mCSD 3.9.0 with CP integrated
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 3.8.0 -> 3.9.0
$ws.Range("B3").Value = "3.9.0"

# Experimental: (blank) -> "false"
# Typing the literal word false directly would auto-convert the cell to a
# Boolean, so build the text via a formula in a scratch cell and paste the
# resulting value back in - this keeps B7 a genuine text cell.
$ws.Range("Z1").Formula = '="false"'
$ws.Range("Z1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# Date: 2022-08-12T09:44:57-05:00 -> 2024-12-02T17:05:26-06:00
$ws.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# Contact rows 10-12 previously all shared the same placeholder text;
# give each its own resolved contact value.
$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws.Range("B11").Value = "null (iti@ihe.net)"
$ws.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction: World -> Global (Whole world)
$ws.Range("B13").Value = "Global (Whole world)"
